# WAHPConfigFile.xlsx update
# - Updated package version for UiPath.UIAutomation.Activities
# - Updated ConfigFile for package uploading in Orchestrator
#
# This script mutates the "Config" and "ConfigOptions" sheets so the
# DownloadsFolder / recipient / base-folder values reflect the new
# deployment target (E:\Bot_Files\RPA FL Renaming\FLOBOT, the
# svc-RCOUIPBOT0005 downloads folder, and the updated CC/To recipient
# lists), and also updates the last-used selection on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Config"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config")

# BaseURL value (row 6) -> new FLOBOT path
$ws1.Range("B6").Value = "E:\Bot_Files\RPA FL Renaming\FLOBOT\"

# SuffixToday-adjacent MasterFolder value (row 8) -> new FLOBOT path
# (also normalize its style back to the plain wrap style used elsewhere,
# matching the style used by the surrounding "Value" column cells)
$ws1.Range("B7").Copy() | Out-Null
$ws1.Range("B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws1.Range("B8").Value = "E:\Bot_Files\RPA FL Renaming\FLOBOT"

# DownloadsFolder value (row 12) -> new service-account downloads path
$ws1.Range("B12").Value = "C:\Users\svc-RCOUIPBOT0005\Downloads\"

# Restore selection / scroll position used in the saved workbook
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B6:B12").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "ConfigOptions"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ConfigOptions")

# RecipientTo / RecipientCC for the first (WAHP) block
$ws2.Range("B11").Value = "doris.cuaresma@lexisnexisrisk.com; shiela.agravante@lexisnexisrisk.com"
$ws2.Range("B12").Value = "agnes.sara@lexisnexisrisk.com; gerard.mancenido@lexisnexisrisk.com; lester.rollan@lexisnexisrisk.com; dindee.galindo@lexisnexisrisk.com; jesriel.tolentino@lexisnexisrisk.com; paul.fabro@lexisnexisrisk.com; judy.cotaoco@lexisnexisrisk.com"
$ws2.Rows.Item(12).RowHeight = 58

# RecipientTo / RecipientCC for the second (FLOBOT testing) block
$ws2.Range("B25").Value = "doris.cuaresma@lexisnexisrisk.com; shiela.agravante@lexisnexisrisk.com"
$ws2.Range("B26").Value = "agnes.sara@lexisnexisrisk.com; gerard.mancenido@lexisnexisrisk.com; lester.rollan@lexisnexisrisk.com; dindee.galindo@lexisnexisrisk.com; jesriel.tolentino@lexisnexisrisk.com; paul.fabro@lexisnexisrisk.com; judy.cotaoco@lexisnexisrisk.com"

# RecipientTo / RecipientCC for the third (virtual desktop testing) block
$ws2.Range("B39").Value = "doris.cuaresma@lexisnexisrisk.com; shiela.agravante@lexisnexisrisk.com"
$ws2.Range("B40").Value = "agnes.sara@lexisnexisrisk.com; gerard.mancenido@lexisnexisrisk.com; lester.rollan@lexisnexisrisk.com; dindee.galindo@lexisnexisrisk.com; jesriel.tolentino@lexisnexisrisk.com; paul.fabro@lexisnexisrisk.com; judy.cotaoco@lexisnexisrisk.com"

# Restore selection / scroll position used in the saved workbook
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("B35:B41").Select() | Out-Null

# Leave "Config" as the active/selected sheet, as in the saved workbook
$ws1.Activate()
